$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns C, D, E, F, G, K across rows 2-25
# Row order per row: C, D, E, F, G, K
$data = @{
    2 = @(0.04837695099705286, 0.1654233272190879, 0.4065214141587177, 3.263238663950631, 0.002479414180627874, 2.08322645173314)
    3 = @(0.0429487227180374, 0.1556740527298217, 0.3541088049937287, 3.05905824376805, 0.002487328911097215, 1.900019779856621)
    4 = @(0.03964073651280842, 0.1497338257248089, 0.3221122585896552, 2.935197352526131, 0.002492423977776161, 1.789369609338735)
    5 = @(0.03829870345259678, 0.1473237035714874, 0.3091156347294799, 2.885086278962575, 0.002494559739637043, 1.74473124387606)
    6 = @(0.03807621408034834, 0.146924114931636, 0.3069599818509943, 2.876786859891695, 0.002494917982069, 1.737346116877177)
    7 = @(0.039622613465923, 0.1497012804508557, 0.321936816207554, 2.934520086790769, 0.002492452540097843, 1.788765781767665)
    8 = @(0.04649994386571166, 0.1620516711690101, 0.3884085248336646, 3.19251500449235, 0.002482094499877864, 2.01966880485827)
    9 = @(0.06019800091988259, 0.186679564883093, 0.5204300962234498, 3.711134600427073, 0.002463636800355556, 2.487555486019744)
    10 = @(0.0704110083863867, 0.2050870733451973, 0.6187393870277447, 4.100970710706349, 0.00245118758800939, 2.841243672064536)
    11 = @(0.07509405223717636, 0.2135433964136837, 0.6638141786506679, 4.280467923527056, 0.002445761417285965, 3.004469851095109)
    12 = @(0.07687309338791692, 0.2167586280764624, 0.6809390667023791, 4.348768513992241, 0.002443740431504835, 3.066628143610103)
    13 = @(0.07648968717735727, 0.2160655739076844, 0.6772483383592203, 4.33404383374392, 0.002444174189086037, 3.05322555364063)
    14 = @(0.07524029932467613, 0.2138076483067266, 0.6652218998418675, 4.286080355473644, 0.002445594474094179, 3.009576600278706)
    15 = @(0.07447576237379394, 0.2124263310318781, 0.657862815079838, 4.25674472664781, 0.002446468830947972, 2.982886073818236)
    16 = @(0.07010573574208934, 0.2045361911970076, 0.61580119785188, 4.089285154558979, 0.002451546945317911, 2.830624459461831)
    17 = @(0.06743460575842164, 0.1997178056749647, 0.5900919421791855, 3.987120052992509, 0.002454722703547811, 2.737823043102253)
    18 = @(0.06590171314724103, 0.1969541093306759, 0.5753376006587416, 3.928559456994094, 0.002456571644817052, 2.684664796602135)
    19 = @(0.06538328795100767, 0.1960196599573578, 0.5703475353533207, 3.908765995372931, 0.002457201508925929, 2.666703502653661)
    20 = @(0.06771859046922657, 0.2002299248807162, 0.5928252895634074, 3.9979746568265, 0.002454382329846758, 2.747679189165126)
    21 = @(0.07560711807390419, 0.2144704937817039, 0.6687527929912136, 4.30015931643743, 0.002445176387455403, 3.022387804933203)
    22 = @(0.0807959644309193, 0.2238538720026497, 0.7187050005081375, 4.49958036885107, 0.002439356569602128, 3.203960910546982)
    23 = @(0.07802342343930491, 0.2188384249788271, 0.6920126797773349, 4.392963167848393, 0.002442444808077883, 3.106861365772886)
    24 = @(0.0675901923613651, 0.1999983757579002, 0.5915894618350279, 3.993066745237513, 0.002454536140761856, 2.743222620593258)
    25 = @(0.0564676171431131, 0.1799667363336539, 0.4845059269482448, 3.569366308597097, 0.002468433510945671, 2.359290125037148)
}

$columns = @("C", "D", "E", "F", "G", "K")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range("$($columns[$i])$row").Value = $values[$i]
    }
}